# Release v0.1.0-beta: Fix validation errors and update canonical URL
#
# Applies:
#   1. Metadata sheet: Version, Status, Experimental, Date, Description updates
#   2. "Include #0" sheet: replace the old Operation/concept/is-a table with a
#      Concept/Description table (3 SNOMED concepts) + System URI row
#   3. New "Include #1" sheet (duplicate of "Include #0" layout) with a single
#      local "risk-falls" concept + System URI row pointing at the local CodeSystem

# Helper: force a numeric-looking / boolean-looking string (e.g. "false",
# "129839007") to be written as literal text instead of being auto-coerced
# by Excel's smart input parsing. We do this by writing a formula that
# evaluates to the text value, then converting the formula to its value
# in place (Copy + PasteSpecial values-only) so the final cell is a plain
# text cell using the ORIGINAL cell style (no quote-prefix style is added).
function Set-TextValue($range, [string]$text) {
    $escaped = $text -replace '"', '""'
    $range.Formula = '=TEXT("' + $escaped + '","@")'
    $range.Copy()
    $range.PasteSpecial(-4163)
}

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. Metadata sheet
# ---------------------------------------------------------------------
$meta = $wb.Worksheets.Item("Metadata")

$meta.Range("B3").Value = "0.1.0"
$meta.Range("B6").Value = "draft"
Set-TextValue $meta.Range("B7") "false"
$meta.Range("B8").Value = "2025-12-26T14:13:58+00:00"
$meta.Range("B11").Value = "Value set for nursing problems and diagnoses"

# ---------------------------------------------------------------------
# 2. Duplicate "Include #0" BEFORE editing it, so the new "Include #1"
#    sheet inherits the exact same column widths / styles the original
#    "Include #0" sheet had (A:30.703125, B:50.703125).
# ---------------------------------------------------------------------
$inc0 = $wb.Worksheets.Item("Include #0")
$inc0.Copy($null, $inc0)
$inc1 = $wb.Worksheets.Item($wb.Worksheets.Count)
$inc1.Name = "Include #1"

# ---------------------------------------------------------------------
# 3. Rewrite "Include #0": drop column C, add 2 rows, new concept table
# ---------------------------------------------------------------------
$inc0.Columns("C").Delete()

# Extend formatting (style) of existing data row 4 down into new rows 5-6
$inc0.Range("A4:B4").Copy()
$inc0.Range("A5:B6").PasteSpecial(-4122)

$inc0.Range("A1").Value = "Concept"
$inc0.Range("B1").Value = "Description"

Set-TextValue $inc0.Range("A2") "129839007"
$inc0.Range("B2").Value = "At risk for falls"

Set-TextValue $inc0.Range("A3") "300893006"
$inc0.Range("B3").Value = "Nutritional finding"

Set-TextValue $inc0.Range("A4") "22253000"
$inc0.Range("B4").Value = "Pain"

$inc0.Range("A5").ClearContents()
$inc0.Range("B5").ClearContents()

$inc0.Range("A6").Value = "System URI"
$inc0.Range("B6").Value = "http://snomed.info/sct"

# ---------------------------------------------------------------------
# 4. Rewrite "Include #1" (the copy): drop column C, new concept table
# ---------------------------------------------------------------------
$inc1.Columns("C").Delete()

$inc1.Range("A1").Value = "Concept"
$inc1.Range("B1").Value = "Description"

Set-TextValue $inc1.Range("A2") "risk-falls"
$inc1.Range("B2").Value = "Risk of falls"

$inc1.Range("A3").ClearContents()
$inc1.Range("B3").ClearContents()

$inc1.Range("A4").Value = "System URI"
$inc1.Range("B4").Value = "https://clinyqai.github.io/open-nursing-core-ig/CodeSystem/onc-observation-codes"

# ---------------------------------------------------------------------
# 5. Keep "Metadata" the active/selected sheet (matches activeTab="0")
# ---------------------------------------------------------------------
$meta.Activate()
